$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Insert three new rows above the current row 3 (shifts the existing
#    walk-in records, old rows 3-8, down to rows 6-11).
# ---------------------------------------------------------------------
$ws.Rows("3:5").Insert()

# Newly inserted rows inherit the formatting of the row above (the
# header row), so reset them back to the default/"Normal" style that the
# rest of the data rows use before writing any values into them.
$ws.Range("A3:K5").Style = "Normal"
$ws.Range("A3:A5").NumberFormat = "0"
$ws.Range("D3:D5").NumberFormat = "0"

# ---------------------------------------------------------------------
# 2. Populate the three newly inserted rows (new walk-in entries #1-#3).
# ---------------------------------------------------------------------
$newTop = @(
    @(1, "10-12-2025", "Rahul",   9562201952, "24-12-2025", "VISHNU N", "Loss", "PRODUCT",                   "PRODUCT NOT AVAILABLE", "-", "kurtha collection not available"),
    @(2, "11-12-2025", "shafnas", 8136863051, "20-12-2025", "ARJUN P",  "Loss", "CUSTOMER INTERNAL ISSUES",  "FAMILY DISAPPROVEL",    "-", "Tommorow coming"),
    @(3, "12-12-2025", "nihal",   9037517955, "24-01-2026", "ARJUN P",  "Loss", "CUSTOMER INTERNAL ISSUES",  "FAMILY DISAPPROVEL",    "-", "next month function afeter coming and booking")
)

$r = 3
foreach ($row in $newTop) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
    $ws.Cells.Item($r, 9).Value = $row[8]
    $ws.Cells.Item($r, 10).Value = $row[9]
    $ws.Cells.Item($r, 11).Value = $row[10]
    $r = $r + 1
}

# ---------------------------------------------------------------------
# 3. Renumber the "#" column for the original records that were pushed
#    down to rows 6-11 (they used to be 1-6, now they continue the
#    sequence as 4-9).
# ---------------------------------------------------------------------
$ws.Range("A6").Value = 4
$ws.Range("A7").Value = 5
$ws.Range("A8").Value = 6
$ws.Range("A9").Value = 7
$ws.Range("A10").Value = 8
$ws.Range("A11").Value = 9

# ---------------------------------------------------------------------
# 4. Append three brand-new walk-in records after the existing data
#    (new rows 12-14, continuing the "#" sequence as 10-12).
# ---------------------------------------------------------------------
$ws.Range("A12:A14").NumberFormat = "0"
$ws.Range("D12:D14").NumberFormat = "0"

$newBottom = @(
    @(10, "19-12-2025", "ABINAV", 9744047563, "04-01-2026", "VISHNU N", "Loss", "ENQUIRY", "ENQUIRY WITHOUT BRIDE/FAMILY",  "-", "just checking"),
    @(11, "19-12-2025", "Ajmal",  8943423460, "26-12-2025", "VISHNU N", "Loss", "PRODUCT", "REQUIRED DESIGN NOT AVAILABLE", "-", "will cheke with family and call"),
    @(12, "20-12-2025", "000",    9496101363, "03-01-2026", "VISHNU N", "Loss", "PRODUCT", "REQUIRED MODEL NOT AVAILABLE",  "-", "wantbmore collection")
)

$r = 12
foreach ($row in $newBottom) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
    $ws.Cells.Item($r, 9).Value = $row[8]
    $ws.Cells.Item($r, 10).Value = $row[9]
    $ws.Cells.Item($r, 11).Value = $row[10]
    $r = $r + 1
}

# ---------------------------------------------------------------------
# 5. Widen columns I ("Sub Category") and K ("Remarks").
#    (Input values chosen so the engine's pixel-snapped ColumnWidth
#    setter lands on the closest achievable value to the target
#    39.150000000000006 / 54 stored widths.)
# ---------------------------------------------------------------------
$ws.Columns("I").ColumnWidth = 38.25
$ws.Columns("K").ColumnWidth = 53.15
